$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 51503.9
$ws.Range("I70").Value = 143768.28
$ws.Range("J70").Value = 1823.0769
$ws.Range("K70").Value = 431304.84
$ws.Range("L70").Value = 5469.2307
$ws.Range("M70").Value = -431034.84
$ws.Range("N70").Value = -6009.2307
$ws.Range("H73").Value = 51503.9
$ws.Range("I73").Value = 143768.28
$ws.Range("J73").Value = 1823.0769
$ws.Range("K73").Value = 431304.84
$ws.Range("L73").Value = 5469.2307
$ws.Range("M73").Value = -430368.84
$ws.Range("N73").Value = -7341.2307
$ws.Range("H116").Value = 1700.4166
$ws.Range("J116").Value = 1800
$ws.Range("L116").Value = 1800
$ws.Range("N116").Value = -8684
$ws.Range("H129").Value = 19389.76
$ws.Range("J129").Value = 25355.023
$ws.Range("L129").Value = 76065.069
$ws.Range("N129").Value = -86065.069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1055.2258
$ws.Range("I2").Value = 911.1111
$ws.Range("J2").Value = 1254.7693
$ws.Range("K2").Value = 911.1111
$ws.Range("L2").Value = 1254.7693
$ws.Range("M2").Value = -798.1111
$ws.Range("N2").Value = -1480.7693
$ws.Range("H45").Value = 1567.4286
$ws.Range("I45").Value = 1558.3334
$ws.Range("J45").Value = 1574.25
$ws.Range("K45").Value = 1558.3334
$ws.Range("L45").Value = 1574.25
$ws.Range("M45").Value = -1181.3334
$ws.Range("N45").Value = -2328.25
$ws.Range("H74").Value = 874.4545000000001
$ws.Range("I74").Value = 811.9
$ws.Range("K74").Value = 811.9
$ws.Range("M74").Value = 62.10000000000002
$ws.Range("H77").Value = 874.4545000000001
$ws.Range("I77").Value = 811.9
$ws.Range("K77").Value = 4059.5
$ws.Range("M77").Value = 308.5
$ws.Range("H116").Value = 1055.2258
$ws.Range("I116").Value = 911.1111
$ws.Range("J116").Value = 1254.7693
$ws.Range("K116").Value = 911.1111
$ws.Range("L116").Value = 1254.7693
$ws.Range("M116").Value = 1382.8889
$ws.Range("N116").Value = -5842.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1055.2258
$ws.Range("I3").Value = 911.1111
$ws.Range("J3").Value = 1254.7693
$ws.Range("K3").Value = 911.1111
$ws.Range("L3").Value = 1254.7693
$ws.Range("M3").Value = -797.1111
$ws.Range("N3").Value = -1482.7693
$ws.Range("H20").Value = 3682.7827
$ws.Range("I20").Value = 3994.2778
$ws.Range("J20").Value = 2561.4
$ws.Range("K20").Value = 3994.2778
$ws.Range("L20").Value = 2561.4
$ws.Range("M20").Value = -3747.2778
$ws.Range("N20").Value = -3055.4
$ws.Range("H86").Value = 2444.4443
$ws.Range("I86").Value = 1428.5714
$ws.Range("K86").Value = 1428.5714
$ws.Range("M86").Value = -305.5714
$ws.Range("H89").Value = 2444.4443
$ws.Range("I89").Value = 1428.5714
$ws.Range("K89").Value = 7142.857
$ws.Range("M89").Value = -1526.857
$ws.Range("H107").Value = 1081.375
$ws.Range("I107").Value = 720.3333
$ws.Range("K107").Value = 720.3333
$ws.Range("M107").Value = 1199.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3706089.8
$ws.Range("I31").Value = 2025.4706
$ws.Range("K31").Value = 2025.4706
$ws.Range("M31").Value = -1730.4706
$ws.Range("H34").Value = 3706089.8
$ws.Range("I34").Value = 2025.4706
$ws.Range("K34").Value = 2025.4706
$ws.Range("M34").Value = -1823.4706
$ws.Range("H99").Value = 2328.9565
$ws.Range("I99").Value = 2438.75
$ws.Range("J99").Value = 2270.4
$ws.Range("K99").Value = 2438.75
$ws.Range("L99").Value = 2270.4
$ws.Range("M99").Value = -940.75
$ws.Range("N99").Value = -5266.4
$ws.Range("H122").Value = 806.9524
$ws.Range("I122").Value = 802.5
$ws.Range("J122").Value = 812.8889
$ws.Range("K122").Value = 2407.5
$ws.Range("L122").Value = 2438.6667
$ws.Range("M122").Value = 42.5
$ws.Range("N122").Value = -7338.6667
$ws.Range("H126").Value = 2328.9565
$ws.Range("I126").Value = 2438.75
$ws.Range("J126").Value = 2270.4
$ws.Range("K126").Value = 7316.25
$ws.Range("L126").Value = 6811.200000000001
$ws.Range("M126").Value = -4846.25
$ws.Range("N126").Value = -11751.2
$ws.Range("H132").Value = 2177.4119
$ws.Range("I132").Value = 1987.5555
$ws.Range("J132").Value = 2391
$ws.Range("K132").Value = 5962.666499999999
$ws.Range("L132").Value = 7173
$ws.Range("M132").Value = -3432.666499999999
$ws.Range("N132").Value = -12233

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1456.5
$ws.Range("I5").Value = 295.72726
$ws.Range("J5").Value = 2207.5881
$ws.Range("K5").Value = 887.18178
$ws.Range("L5").Value = 6622.7643
$ws.Range("M5").Value = -775.18178
$ws.Range("N5").Value = -6846.7643
$ws.Range("H131").Value = 2454.6177
$ws.Range("I131").Value = 18808.334
$ws.Range("J131").Value = 872
$ws.Range("K131").Value = 56425.00199999999
$ws.Range("L131").Value = 2616
$ws.Range("M131").Value = -51385.00199999999
$ws.Range("N131").Value = -12696
$ws.Range("H135").Value = 1456.5
$ws.Range("I135").Value = 295.72726
$ws.Range("J135").Value = 2207.5881
$ws.Range("K135").Value = 2661.54534
$ws.Range("L135").Value = 19868.2929
$ws.Range("M135").Value = -126.5453400000001
$ws.Range("N135").Value = -24938.2929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5129.4
$ws.Range("I80").Value = 4777.6
$ws.Range("J80").Value = 5481.2
$ws.Range("K80").Value = 4777.6
$ws.Range("L80").Value = 5481.2
$ws.Range("M80").Value = -3779.6
$ws.Range("N80").Value = -7477.2
$ws.Range("H83").Value = 5129.4
$ws.Range("I83").Value = 4777.6
$ws.Range("J83").Value = 5481.2
$ws.Range("K83").Value = 23888
$ws.Range("L83").Value = 27406
$ws.Range("M83").Value = -18896
$ws.Range("N83").Value = -37390
$ws.Range("H126").Value = 3119.2307
$ws.Range("I126").Value = 4520
$ws.Range("J126").Value = 2243.75
$ws.Range("K126").Value = 13560
$ws.Range("L126").Value = 6731.25
$ws.Range("M126").Value = -11090
$ws.Range("N126").Value = -11671.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 258
$ws.Range("I55").Value = 264.5
$ws.Range("K55").Value = 264.5
$ws.Range("M55").Value = -91.5
$ws.Range("H132").Value = 2741.8333
$ws.Range("I132").Value = 2934
$ws.Range("J132").Value = 2549.6667
$ws.Range("K132").Value = 8802
$ws.Range("L132").Value = 7649.000100000001
$ws.Range("M132").Value = -6272
$ws.Range("N132").Value = -12709.0001
$ws.Range("H133").Value = 24605.2
$ws.Range("J133").Value = 24605.2
$ws.Range("L133").Value = 24605.2
$ws.Range("N133").Value = -29665.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3321.6365
$ws.Range("I122").Value = 3800.8
$ws.Range("J122").Value = 2922.3333
$ws.Range("K122").Value = 11402.4
$ws.Range("L122").Value = 8766.999899999999
$ws.Range("M122").Value = -8952.400000000001
$ws.Range("N122").Value = -13666.9999
$ws.Range("H132").Value = 2722.2222
$ws.Range("I132").Value = 1501
$ws.Range("J132").Value = 3699.2
$ws.Range("K132").Value = 4503
$ws.Range("L132").Value = 11097.6
$ws.Range("M132").Value = -1973
$ws.Range("N132").Value = -16157.6
